# Auto-generated edit script for Línea 141 schedule refresh (scrape update 12:46:07)
$wb = $excel.ActiveWorkbook

# ---------- Sheet: LP1912 ----------
$ws = $wb.Worksheets.Item('LP1912')

$ws.Cells.Item(2, 1).Value = 'Última actualización: 12:46:07'
$ws.Cells.Item(3, 1).Value = 'Total filas: 166'
$ws.Cells.Item(45, 1).Value = '08:45:31'
$ws.Cells.Item(45, 3).Value = '215C_EL PATO'
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(46, 1).Value = '07:56:02'
$ws.Cells.Item(46, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(46, 4).Value = 49
$ws.Cells.Item(66, 1).Value = '08:28:52'
$ws.Cells.Item(66, 3).Value = '10_OLMOS'
$ws.Cells.Item(66, 4).Value = 60
$ws.Cells.Item(67, 1).Value = '08:11:18'
$ws.Cells.Item(67, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(67, 4).Value = 77
$ws.Cells.Item(118, 1).Value = '10:36:50'
$ws.Cells.Item(118, 3).Value = '225_GOMEZ'
$ws.Cells.Item(118, 4).Value = 76
$ws.Cells.Item(119, 1).Value = '11:33:52'
$ws.Cells.Item(119, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(119, 4).Value = 19
$ws.Cells.Item(139, 1).Value = '11:53:44'
$ws.Cells.Item(139, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(139, 4).Value = 43
$ws.Cells.Item(140, 1).Value = '10:49:38'
$ws.Cells.Item(140, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(140, 4).Value = 107
$ws.Cells.Item(152, 1).Value = '12:46:07'
$ws.Cells.Item(152, 2).Value = '13:19'
$ws.Cells.Item(152, 3).Value = '15_ABASTO'
$ws.Cells.Item(152, 4).Value = 33
$ws.Cells.Item(153, 1).Value = '11:53:44'
$ws.Cells.Item(153, 2).Value = '13:21'
$ws.Cells.Item(153, 4).Value = 88
$ws.Cells.Item(154, 1).Value = '12:46:07'
$ws.Cells.Item(154, 2).Value = '13:22'
$ws.Cells.Item(154, 4).Value = 36
$ws.Cells.Item(155, 1).Value = '12:33:02'
$ws.Cells.Item(155, 2).Value = '13:23'
$ws.Cells.Item(155, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(155, 4).Value = 50
$ws.Cells.Item(156, 2).Value = '13:24'
$ws.Cells.Item(156, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(156, 4).Value = 73
$ws.Cells.Item(157, 1).Value = '11:33:52'
$ws.Cells.Item(157, 2).Value = '13:25'
$ws.Cells.Item(157, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(157, 4).Value = 112
$ws.Cells.Item(158, 1).Value = '12:11:21'
$ws.Cells.Item(158, 2).Value = '13:32'
$ws.Cells.Item(158, 3).Value = '14_ABASTO'
$ws.Cells.Item(158, 4).Value = 81
$ws.Cells.Item(159, 1).Value = '11:53:44'
$ws.Cells.Item(159, 2).Value = '13:32'
$ws.Cells.Item(159, 3).Value = '215A_EL PATO'
$ws.Cells.Item(159, 4).Value = 99
$ws.Cells.Item(160, 1).Value = '11:46:32'
$ws.Cells.Item(160, 2).Value = '13:33'
$ws.Cells.Item(160, 3).Value = '215A_EL PATO'
$ws.Cells.Item(160, 4).Value = 107
$ws.Cells.Item(161, 2).Value = '13:33'
$ws.Cells.Item(161, 3).Value = '14_ABASTO'
$ws.Cells.Item(161, 4).Value = 60
$ws.Cells.Item(162, 1).Value = '11:53:44'
$ws.Cells.Item(162, 2).Value = '13:47'
$ws.Cells.Item(162, 3).Value = '225_GOMEZ'
$ws.Cells.Item(162, 4).Value = 114
$ws.Cells.Item(163, 2).Value = '13:54'
$ws.Cells.Item(163, 3).Value = '15_ABASTO'
$ws.Cells.Item(163, 4).Value = 81
$ws.Cells.Item(164, 2).Value = '14:02'
$ws.Cells.Item(164, 3).Value = '10_OLMOS'
$ws.Cells.Item(164, 4).Value = 89
$ws.Cells.Item(165, 1).Value = '12:46:07'
$ws.Cells.Item(165, 2).Value = '14:02'
$ws.Cells.Item(165, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(165, 4).Value = 76
$ws.Cells.Item(166, 1).Value = '12:46:07'
$ws.Cells.Item(166, 2).Value = '14:08'
$ws.Cells.Item(166, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(166, 4).Value = 82
$ws.Cells.Item(166, 5).Value = 'LP1912'
$ws.Cells.Item(167, 1).Value = '12:33:02'
$ws.Cells.Item(167, 2).Value = '14:17'
$ws.Cells.Item(167, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(167, 4).Value = 104
$ws.Cells.Item(167, 5).Value = 'LP1912'
$ws.Cells.Item(168, 1).Value = '12:33:02'
$ws.Cells.Item(168, 2).Value = '14:18'
$ws.Cells.Item(168, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(168, 4).Value = 105
$ws.Cells.Item(168, 5).Value = 'LP1912'
$ws.Cells.Item(169, 1).Value = '12:33:02'
$ws.Cells.Item(169, 2).Value = '14:32'
$ws.Cells.Item(169, 3).Value = '14X44_ABASTO'
$ws.Cells.Item(169, 4).Value = 119
$ws.Cells.Item(169, 5).Value = 'LP1912'
$ws.Cells.Item(170, 1).Value = '12:46:07'
$ws.Cells.Item(170, 2).Value = '14:34'
$ws.Cells.Item(170, 3).Value = '215C_EL PATO'
$ws.Cells.Item(170, 4).Value = 108
$ws.Cells.Item(170, 5).Value = 'LP1912'
$ws.Cells.Item(171, 1).Value = '12:46:07'
$ws.Cells.Item(171, 2).Value = '14:39'
$ws.Cells.Item(171, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(171, 4).Value = 113
$ws.Cells.Item(171, 5).Value = 'LP1912'

# ---------- Sheet: LP1912-215 ----------
$ws = $wb.Worksheets.Item('LP1912-215')

$ws.Cells.Item(2, 1).Value = 'Última actualización: 12:46:07'
$ws.Cells.Item(3, 1).Value = 'Total filas: 29'
$ws.Cells.Item(34, 1).Value = '12:46:07'
$ws.Cells.Item(34, 2).Value = '14:34'
$ws.Cells.Item(34, 3).Value = '215C_EL PATO'
$ws.Cells.Item(34, 4).Value = 108
$ws.Cells.Item(34, 5).Value = 'LP1912'

# ---------- Sheet: 6203-6173 ----------
$ws = $wb.Worksheets.Item('6203-6173')

$ws.Cells.Item(2, 1).Value = 'Última actualización: 12:46:07'
$ws.Cells.Item(3, 1).Value = 'Total filas: 25'
$ws.Cells.Item(30, 1).Value = '12:46:07'
$ws.Cells.Item(30, 2).Value = '14:27'
$ws.Cells.Item(30, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(30, 4).Value = 101
$ws.Cells.Item(30, 5).Value = 'L6203'

